$d = $word.ActiveDocument

# --- 1. AWO Rössing-Barnten: add the member count "(90)" after the name ---
# Split into two runs: "AWO Rössing-Barnten " and "(90)", matching how the
# entry looks for other clubs in the list (e.g. "CDU Ortsverband Rössing (28)").
$rng = $d.Content
$found = $rng.Find.Execute("AWO Rössing-Barnten", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $origStart = $rng.Start
    $rng.Collapse(0)
    $rng.InsertAfter(" ")
    $afterSpaceEnd = $rng.End
    $rng.Collapse(0)
    $rng.InsertAfter("(90)")
    # Toggling a character property back to its original value forces the
    # engine to keep this newly-inserted text as its own run instead of
    # silently re-merging it with its neighbour.
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0

    $combinedRng = $d.Range($origStart, $afterSpaceEnd)
    $combinedRng.Font.Bold = 1
    $combinedRng.Font.Bold = 0

    Write-Host "AWO Rössing-Barnten: inserted (90)"
} else {
    Write-Host "AWO Rössing-Barnten: NOT FOUND"
}

# --- 2. Landfrauenverein (132): the member count used to be split across
# three runs ("Landfrauenverein (", "132", ")"); normalize back to one run. ---
$found2a = $d.Content.Find.Execute("Landfrauenverein (132)", $true, $false, $false, $false, $false, $true, 1, $false, "Landfrauenverein (132)", 2)
if ($found2a) {
    $rng2 = $d.Content
    $found2b = $rng2.Find.Execute("Landfrauenverein (132)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2b) {
        $rng2.Font.Bold = 1
        $rng2.Font.Bold = 0
    }
    Write-Host "Landfrauenverein (132): runs normalized"
} else {
    Write-Host "Landfrauenverein (132): NOT FOUND"
}

# --- 3. CDU Ortsverband Rössing (28): likewise merge the two existing runs
# ("CDU Ortsverband Rössing ", "(28)") back into a single run. ---
$found3a = $d.Content.Find.Execute("CDU Ortsverband Rössing (28)", $true, $false, $false, $false, $false, $true, 1, $false, "CDU Ortsverband Rössing (28)", 2)
if ($found3a) {
    $rng3 = $d.Content
    $found3b = $rng3.Find.Execute("CDU Ortsverband Rössing (28)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found3b) {
        $rng3.Font.Bold = 1
        $rng3.Font.Bold = 0
    }
    Write-Host "CDU Ortsverband Rössing (28): runs normalized"
} else {
    Write-Host "CDU Ortsverband Rössing (28): NOT FOUND"
}
